$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to text,
# matching the source data which stores these as formatted strings
# (e.g. "1.00", thousands-dot separated, etc.) rather than numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '66.930.68'
$ws.Range("E2").Value = '  -4.20%  '
$ws.Range("D3").Value = '3.518.27'
$ws.Range("E3").Value = '  -4.67%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws.Range("D5") '607.97'
$ws.Range("E5").Value = '  -6.29%  '
Set-TextValue $ws.Range("D6") '152.09'
$ws.Range("E6").Value = '  -5.53%  '
$ws.Range("D7").Value = '3.516.07'
$ws.Range("E7").Value = '  -4.61%  '
$ws.Range("E8").Value = '  +0.09%  '
Set-TextValue $ws.Range("D9") '0.483'
$ws.Range("E9").Value = '  -3.99%  '
$ws.Range("E10").Value = '  -4.68%  '
$ws.Range("E11").Value = '  -4.16%  '
Set-TextValue $ws.Range("D12") '0.427'
$ws.Range("E12").Value = '  -3.96%  '
$ws.Range("E13").Value = '  -5.38%  '
$ws.Range("D14").Value = '4.116.91'
$ws.Range("E14").Value = '  -4.57%  '
Set-TextValue $ws.Range("D15") '31.61'
$ws.Range("E15").Value = '  -3.62%  '
$ws.Range("D16").Value = '3.532.23'
$ws.Range("E16").Value = '  -4.26%  '
$ws.Range("D17").Value = '66.974.16'
$ws.Range("E17").Value = '  -4.14%  '
$ws.Range("E18").Value = '  +0.26%  '
Set-TextValue $ws.Range("D19") '6.32'
$ws.Range("E19").Value = '  -3.01%  '
$ws.Range("E20").Value = '  -4.97%  '
Set-TextValue $ws.Range("D21") '445.77'
$ws.Range("E21").Value = '  -5.41%  '
Set-TextValue $ws.Range("D22") '8.98'
$ws.Range("E22").Value = '  -14.35%  '
Set-TextValue $ws.Range("D23") '0.629'
$ws.Range("E23").Value = '  -3.27%  '
Set-TextValue $ws.Range("D24") '77.28'
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '3.661.62'
$ws.Range("E26").Value = '  -4.58%  '
$ws.Range("E27").Value = '  -3.25%  '
Set-TextValue $ws.Range("D28") '10.16'
$ws.Range("E28").Value = '  -7.06%  '
$ws.Range("E29").Value = '  -10.71%  '
$ws.Range("E30").Value = '  -5.01%  '
Set-TextValue $ws.Range("D31") '1.61'
$ws.Range("E31").Value = '  -6.31%  '
Set-TextValue $ws.Range("D32") '1.00'
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("E33").Value = '  -1.90%  '
Set-TextValue $ws.Range("D34") '25.66'
$ws.Range("E34").Value = '  -4.43%  '
Set-TextValue $ws.Range("D35") '6.17'
$ws.Range("E35").Value = '  -5.88%  '
Set-TextValue $ws.Range("D36") '1.86'
$ws.Range("E36").Value = '  -7.58%  '
$ws.Range("D37").Value = '3.512.04'
$ws.Range("E37").Value = '  -4.79%  '
Set-TextValue $ws.Range("D38") '8.01'
$ws.Range("E38").Value = '  -5.48%  '
Set-TextValue $ws.Range("D40") '1.00'
$ws.Range("E40").Value = '  -0.04%  '
Set-TextValue $ws.Range("D41") '172.86'
$ws.Range("E41").Value = '  -3.59%  '
Set-TextValue $ws.Range("D42") '2.15'
$ws.Range("E42").Value = '  -4.63%  '
Set-TextValue $ws.Range("D43") '5.54'
$ws.Range("E43").Value = '  -6.21%  '
Set-TextValue $ws.Range("D44") '0.0863'
$ws.Range("E44").Value = '  -4.50%  '
Set-TextValue $ws.Range("D45") '0.889'
$ws.Range("E45").Value = '  -4.30%  '
Set-TextValue $ws.Range("D46") '45.18'
$ws.Range("E46").Value = '  -4.20%  '
Set-TextValue $ws.Range("D47") '27.14'
$ws.Range("E47").Value = '  -6.81%  '
Set-TextValue $ws.Range("D48") '2.53'
$ws.Range("E48").Value = '  -6.88%  '
$ws.Range("E49").Value = '  -1.80%  '
Set-TextValue $ws.Range("D50") '7.56'
$ws.Range("E50").Value = '  -3.64%  '
$ws.Range("E51").Value = '  -5.58%  '
